$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.393.59'
$ws.Range("E2").Value = '  +1.69%  '
$ws.Range("D3").Value = '2.683.01'
$ws.Range("E3").Value = '  +1.88%  '
$ws.Range("E4").Value = '  +0.00%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '600.49'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.04%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '177.33'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -2.05%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '2.682.04'
$ws.Range("E9").Value = '  +1.92%  '
$ws.Range("E10").Value = '  +2.67%  '
$ws.Range("E11").Value = '  +2.16%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.354'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +2.31%  '
$ws.Range("E13").Value = '  +0.33%  '
$ws.Range("D14").Value = '3.173.06'
$ws.Range("E14").Value = '  +3.99%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.0000185'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -2.09%  '
$ws.Range("D16").Value = '72.306.62'
$ws.Range("E16").Value = '  +1.61%  '
$ws.Range("E17").Value = '  -1.26%  '
$ws.Range("D18").Value = '2.687.91'
$ws.Range("E18").Value = '  +2.57%  '
$ws.Range("E19").Value = '  +4.62%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '7.90'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.01%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '371.17'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -2.94%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.17'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.76%  '
$ws.Range("E23").Value = '  +7.38%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '72.21'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("E25").Value = '  -0.02%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '4.35'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -2.93%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.82'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.84%  '
$ws.Range("D28").Value = '2.827.34'
$ws.Range("E28").Value = '  +2.18%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("D30").Value = '0.0₃0941'
$ws.Range("E30").Value = '  -2.50%  '
$ws.Range("E31").Value = '  +0.26%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '511.63'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -5.53%  '
$ws.Range("E33").Value = '  -1.99%  '
$ws.Range("E34").Value = '  -1.09%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +0.01%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '162.67'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.90%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '19.59'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +2.03%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '19.11'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.43%  '
$ws.Range("E39").Value = '  -0.70%  '
$ws.Range("E40").Value = '  -3.77%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.108'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -8.79%  '
$ws.Range("E42").Value = '  -0.05%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '2.57'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -2.96%  '
$ws.Range("E45").Value = '  +0.43%  '
$ws.Range("E46").Value = '  -1.69%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '153.76'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.30%  '
$ws.Range("E48").Value = '  +2.42%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.554'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +3.70%  '
$ws.Range("E50").Value = '  +1.69%  '
$ws.Range("E51").Value = '  +1.56%  '
